$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (C) column date for all existing data rows (2-265)
#    from 2023-09-20 (45189) to 2023-09-21 (45190).
$ws.Range("C2:C265").Value = 45190

# 2. Row 265 picks up an explicit row height (matches default 15pt, but now
#    stored explicitly with customHeight flag).
$ws.Range("A265").EntireRow.RowHeight = 15

# 3. Append a new data row (266) for case "A 44365-2023".
$ws.Range("B266").NumberFormat = "YYYY-MM-DD"
$ws.Range("C266").NumberFormat = "YYYY-MM-DD"

$ws.Range("A266").Value = "A 44365-2023"
$ws.Range("B266").Value = 45188
$ws.Range("C266").Value = 45190
$ws.Range("D266").Value = "SKÅNE LÄN"
$ws.Range("E266").Value = "KLIPPAN"
$ws.Range("G266").Value = 5.1
$ws.Range("H266").Value = 0
$ws.Range("I266").Value = 0
$ws.Range("J266").Value = 0
$ws.Range("K266").Value = 0
$ws.Range("L266").Value = 0
$ws.Range("M266").Value = 0
$ws.Range("N266").Value = 0
$ws.Range("O266").Value = 0
$ws.Range("P266").Value = 0
$ws.Range("Q266").Value = 0

$ws.Range("R266").WrapText = $true
